# "Se agrega para cancelar cupo"
# GIMENEZ ZAIRA had booked Spinning on MARTES 08:00 (row 2). That class/slot
# is now cancelled and she is re-booked into Localizada on MARTES 09:00
# instead, recorded as a new row 3 on her personal sheet. The corresponding
# "cupo" (slot) flags on the Localizada and Spinning schedule sheets are
# flipped to reflect the move.

$wb = $excel.ActiveWorkbook

# --- Localizada: free slot MARTES 09:00 (row 2) gets taken (cupo = 1) ---
$wsLocalizada = $wb.Worksheets.Item("Localizada")
$wsLocalizada.Cells.Item(2, 3).Value = 1

# --- Spinning: slot MARTES 08:00 (row 2) gets cancelled (cupo = 0) ---
$wsSpinning = $wb.Worksheets.Item("Spinning")
$wsSpinning.Cells.Item(2, 3).Value = 0

# --- GIMENEZ ZAIRA: move her booking from Spinning MARTES 08:00 (row 2) to
#     Localizada MARTES 09:00 (new row 3) ---
$wsZaira = $wb.Worksheets.Item("GIMENEZ ZAIRA")
$wsZaira.Cells.Item(2, 1).Value = ""
$wsZaira.Cells.Item(2, 2).Value = ""
$wsZaira.Cells.Item(3, 1).Value = "Localizada"
$wsZaira.Cells.Item(3, 2).Value = "MARTES 09:00"
